$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Break Da Bank Again Respins Free | Slot
#    Game Review") right before the final "Prompt: ..." paragraph. Insert an
#    empty paragraph right after the current second-to-last paragraph (the
#    "Lack of visually impressive graphics" bullet) and then fill that new,
#    still-empty paragraph's range with the target run XML.
$count = $d.Paragraphs.Count
$secondToLastPara = $d.Paragraphs.Item($count - 1)
$secondToLastPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Break Da Bank Again Respins Free | Slot Game Review</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the text of the (now) final paragraph -- swap the image-generation
#    prompt for the meta-description copy, keeping its italic formatting.
$oldText = 'Prompt: Create a cartoon-style feature image for "Break da Bank Again Respins" featuring a happy Maya warrior with glasses. Description: The image should show a Maya warrior wearing glasses and a big smile, holding a bag of gold coins in one hand and a slot machine lever in the other. The background should be filled with colorful banknotes and stacks of gold bars. The Maya warrior should be dressed in traditional clothing, with a feather headdress and colorful patterns on his garment. The overall style of the image should be cartoonish and fun, with bright colors and simple shapes. The image should be eye-catching and convey the excitement and joy of a big win on "Break da Bank Again Respins".'
$newText = 'Read our review of Break da Bank Again Respins slot game. Play now for free and enjoy the unique respin feature and exciting gameplay mechanisms.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
